$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, copying column C (values + formatting) into it
$ws.Columns("C:C").Copy()
$ws.Columns("D:D").Insert()

# Fill in the new "Colonne de(s) maximum(s)" values in column F
$ws.Range("F2").Value = "a"
$ws.Range("F3").Value = "a"
$ws.Range("F4").Value = "a"
$ws.Range("F5").Value = "a"
$ws.Range("F6").Value = "a"
$ws.Range("F7").Value = "af"
$ws.Range("F8").Value = "g"
$ws.Range("F9").Value = "ff"

# Update the active selection
$ws.Range("E9").Select()
